$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sales rows appended below the existing data (rows 13-17)
$newRows = @(
    @(2018, 12, 27, 1, "shirts", 2, 89.90000000000001),
    @(2018, 12, 27, 1, "shirts", 2, 89.90000000000001),
    @(2018, 12, 27, 3, "shoes",  2, 250),
    @(2018, 12, 27, 4, "coats", 1, 350),
    @(2018, 12, 27, 1, "shirts", 1, 89.90000000000001)
)

$startRow = 13
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
}
